$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.674.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.485.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.07%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.17%  "

$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.484.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.27%  "

$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.928.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.513.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.67%  "

$ws.Range("E17").Value = "  +6.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.482.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.85%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.62%  "

$ws.Range("E25").Value = "  +2.35%  "

$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("E27").Value = "  +8.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0831"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.85%  "

$ws.Range("E31").Value = "  +15.23%  "

$ws.Range("E32").Value = "  +7.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "

$ws.Range("E34").Value = "  +10.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.402"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "373.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.89%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("E40").Value = "  +11.92%  "

$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("E42").Value = "  +5.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "151.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.51%  "

$ws.Range("E44").Value = "  +6.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.603"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0968"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0526"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0240"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0228"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.84%  "
